$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 160, shifting existing rows 160:230 down to 161:231
$ws.Rows("160:160").Insert()

# Populate the newly inserted row 160 with the new weekly record
$ws.Cells.Item(160, 1).Value = 10
$ws.Cells.Item(160, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(160, 3).Value = "La Araucanía"
$ws.Cells.Item(160, 4).Value = 44553
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = 100112017
$ws.Cells.Item(160, 7).Value = "Apio"
$ws.Cells.Item(160, 8).Value = "Americana (o)"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 125
$ws.Cells.Item(160, 11).Value = 9000
$ws.Cells.Item(160, 12).Value = 9000
$ws.Cells.Item(160, 13).Value = 9000
$ws.Cells.Item(160, 14).Value = "$/docena de matas"
$ws.Cells.Item(160, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(160, 16).Value = 1500
$ws.Cells.Item(160, 17).Value = 6
$ws.Cells.Item(160, 18).Value = "Hortaliza"
